$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update id column (A) with new numeric values
$ws.Range("A2").Value = 55445201
$ws.Range("A3").Value = 55445202
$ws.Range("A4").Value = 55445203

# Update username column (B) with new usernames
$ws.Range("B2").Value = "testUser55445201"
$ws.Range("B3").Value = "testUser55445202"
$ws.Range("B4").Value = "testUser55445203"
